$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.381.87'
$ws.Range("E2").Value = '  +2.04%  '

$ws.Range("D3").Value = '3.636.75'
$ws.Range("E3").Value = '  +0.74%  '

$ws.Range("E4").Value = '  -0.26%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '197.13'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +7.67%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '579.25'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.38%  '

$ws.Range("D7").Value = '3.630.14'
$ws.Range("E7").Value = '  +1.00%  '

$ws.Range("E8").Value = '  +1.61%  '

$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("E10").Value = '  +1.30%  '

$ws.Range("E11").Value = '  +7.49%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '56.31'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +4.87%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000293'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +16.25%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '10.08'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.92%  '

$ws.Range("D15").Value = '4.217.86'
$ws.Range("E15").Value = '  +0.45%  '

$ws.Range("D16").Value = '3.641.51'
$ws.Range("E16").Value = '  +0.66%  '

$ws.Range("E17").Value = '  +0.70%  '

$ws.Range("E18").Value = '  +3.37%  '

$ws.Range("D19").Value = '68.291.75'
$ws.Range("E19").Value = '  +1.88%  '

$ws.Range("E20").Value = '  +1.83%  '

$ws.Range("E21").Value = '  +3.02%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '403.05'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +2.99%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '13.15'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +28.11%  '

$ws.Range("E24").Value = '  -1.13%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '85.98'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.49%  '

$ws.Range("E26").Value = '  +3.41%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '12.65'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +3.75%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '3.87'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +7.99%  '

$ws.Range("E29").Value = '  +1.12%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '8.16'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +20.75%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '9.18'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.92%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '31.79'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +2.37%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '686.81'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +15.47%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '12.25'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +3.13%  '

$ws.Range("E35").Value = '  +5.81%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '64.73'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.54%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '42.79'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +3.80%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.423'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +13.51%  '

$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("D40").Value = '0.0₃0788'
$ws.Range("E40").Value = '  +7.25%  '

$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.88'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +19.89%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.137'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +4.91%  '

$ws.Range("D43").Value = '3.219.07'
$ws.Range("E43").Value = '  +17.77%  '

$ws.Range("E44").Value = '  +13.44%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.33%  '

$ws.Range("E46").Value = '  +33.11%  '

$ws.Range("E47").Value = '  +2.54%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '8.89'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +7.87%  '

$ws.Range("E49").Value = '  +2.48%  '

$ws.Range("E50").Value = '  +1.80%  '

$ws.Range("E51").Value = '  +3.34%  '

Write-Output "applied 84 cell updates"
